$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: quote_ProposalSetup_263004_TC_06 / 09/06/2022 / Pass
# Row 27: quote_ProposalSetup_263001_TC_03 / 09/06/2022 / Fail
#
# The "Execution Date" column (D) stores dates as plain text strings (not
# real dates) in this workbook. Excel's automatic "looks like a date"
# detection would otherwise convert the text into a date serial number,
# so we briefly force the cells to Text format while entering the value,
# then clear the formatting again so the cells keep the default (no
# explicit style) look used by the rest of the sheet.
$dateCells = $ws.Range("D26:D27")
$dateCells.NumberFormat = "@"

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "quote_ProposalSetup_263004_TC_06"
$ws.Range("D26").Value = "09/06/2022"
$ws.Range("E26").Value = "Pass"

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "quote_ProposalSetup_263001_TC_03"
$ws.Range("D27").Value = "09/06/2022"
$ws.Range("E27").Value = "Fail"

$dateCells.ClearFormats()
